$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (TX001) - update amount/subAmount/vat values
$ws.Range("B2").Value = 81.45
$ws.Range("C2").Value = 65.16
$ws.Range("D2").Value = 16.29

# Row 3 (TX002) - becomes the former "Roaming intäkter - Hubject" data
$ws.Range("B3").Value = 233.65
$ws.Range("C3").Value = 233.65
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "Roaming intäkter - Hubject"
$ws.Range("G3").Value = 30

# Row 4 (TX003) - becomes the former "Plattformsavgift - Monta" data
$ws.Range("B4").Value = -507.7
$ws.Range("C4").Value = -406.16
$ws.Range("D4").Value = -101.54
$ws.Range("E4").Value = 25
$ws.Range("F4").Value = "Plattformsavgift - Monta"
$ws.Range("G4").Value = ""

# Row 5 (TX004) - becomes the former "Övriga kostnader (momsfri)" data
$ws.Range("B5").Value = -20.32
$ws.Range("C5").Value = -20.32
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = "Övriga kostnader (momsfri)"
$ws.Range("G5").Value = ""

# Row 6 (TX005) - removed entirely
$ws.Rows("6").Delete()
